# Append newly-finished books to the reading list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newBooks = @(
    @{
        Title   = "Over Dressed"
        Author  = "Elizabeth Cline"
        Start   = 44204
        Finish  = 44206
        Tags    = "fashion;business;quality;history"
        Type    = "Audio"
        Length  = "7 Hours 57 Mins"
        Rating  = 2
        First   = $true
    },
    @{
        Title   = "Post Corona"
        Author  = "Scott Galloway"
        Start   = 44197
        Finish  = 44207
        Tags    = "business;coronavirus;big tech;adaptation"
        Type    = "Hard Copy"
        Length  = "212 Pages"
        Rating  = 3
        First   = $true
    },
    @{
        Title   = "The Immortal Life of Henrietta Lacks"
        Author  = "Rebecca Skloot"
        Start   = 44206
        Finish  = 44209
        Tags    = "science;ethics;cells;biology;biography;henrietta lacks"
        Type    = "Audio"
        Length  = "12 Hours 13 Mins"
        Rating  = 3
        First   = $true
    },
    @{
        Title   = "Grocery"
        Author  = "Michael Ruhlman"
        Start   = 44209
        Finish  = 44212
        Tags    = "grocery;business;food;health;nutrition"
        Type    = "Audio"
        Length  = "11 Hours 9 Mins"
        Rating  = 4
        First   = $false
    },
    @{
        Title   = "To Pixar and Beyond"
        Author  = "Lawrence Levy"
        Start   = 44207
        Finish  = 44213
        Tags    = "pixar;business;ipo;disney;strategy"
        Type    = "Hard Copy"
        Length  = "248 Pages"
        Rating  = 4
        First   = $false
    }
)

# Grab the date-formatted style from an existing row (row 2) so the new
# Start/Finish Date cells match the rest of the column instead of minting
# a brand-new number format.
$ws.Range("C2:D2").Copy() | Out-Null

$row = 6
foreach ($book in $newBooks) {
    $ws.Cells.Item($row, 1).Value = $book.Title
    $ws.Cells.Item($row, 2).Value = $book.Author

    $ws.Range("C$row`:D$row").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $ws.Cells.Item($row, 3).Value = $book.Start
    $ws.Cells.Item($row, 4).Value = $book.Finish

    $ws.Cells.Item($row, 5).Value = $book.Tags
    $ws.Cells.Item($row, 6).Value = $book.Type
    $ws.Cells.Item($row, 7).Value = $book.Length
    $ws.Cells.Item($row, 8).Value = $book.Rating
    $ws.Cells.Item($row, 9).Value = $book.First
    $row++
}

$excel.CutCopyMode = $false
$ws.Range("A11").Select()
